$d = $word.ActiveDocument
$nl = [char]10

function Replace-Text($search, $replace) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $result = $find.Execute($search, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Host "WARNING: replacement not found for: $search"
    }
    return $result
}

# Line number bumps in the M2DocEvaluator.java / M2DocUtils.java / AbstractTemplatesTestSuite.java stack trace
Replace-Text "M2DocEvaluator.java:1181)" "M2DocEvaluator.java:1204)"
Replace-Text "M2DocEvaluator.java:1216)" "M2DocEvaluator.java:1239)"
Replace-Text "M2DocEvaluator.java:1425)" "M2DocEvaluator.java:1464)"
Replace-Text "M2DocEvaluator.java:287)" "M2DocEvaluator.java:296)"
Replace-Text "M2DocEvaluator.java:276)" "M2DocEvaluator.java:281)"
Replace-Text "M2DocUtils.java:694)" "M2DocUtils.java:805)"
Replace-Text "AbstractTemplatesTestSuite.java:480)" "AbstractTemplatesTestSuite.java:511)"
Replace-Text "AbstractTemplatesTestSuite.java:389)" "AbstractTemplatesTestSuite.java:420)"

# Insert a new stack frame line for RunBefores.evaluate right before the
# RunAfters.evaluate line that directly follows a ParentRunner$2.evaluate line
# (there are several ParentRunner$2.evaluate lines; only this one is
# immediately followed by RunAfters.evaluate, making the search unique).
$search = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl + "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
$replace = "ParentRunner`$2.evaluate(ParentRunner.java:268)" + $nl + "`tat org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + $nl + "`tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)"
Replace-Text $search $replace

Write-Host "Done"
